# Apply the data refresh: shift all timestamps forward by 8 days
# (2025-05-22 -> 2025-05-30, i.e. Excel serial 45799 -> 45807) and
# update the Actual Production values for the first day's worth of
# rows to the newly fetched data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 8 days.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = [double]$cell.Value2 + 8
}

# Updated "Actual Production (MW)" values for rows 2-44.
$newValues = @{
    2  = 236
    3  = 246
    4  = 275
    5  = 304
    6  = 337
    7  = 394
    8  = 473
    9  = 571
    10 = 712
    11 = 796
    12 = 1005
    13 = 1048
    14 = 1095
    15 = 1066
    16 = 1055
    17 = 1066
    18 = 1241
    19 = 1284
    20 = 1263
    21 = 1231
    22 = 1237
    23 = 1233
    24 = 1242
    25 = 1226
    26 = 1206
    27 = 1188
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $newValues[$row]
}
